$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(0.3983698899088211,0.07897745812434209,0.6180439278381868,0.2369271818513674,0,0.002500857343096379,0,0.9477496427068672,0.1120773264191222,0.4589630088180456,0,0.3430080280805754,2.140150110315426,4.855211653973072),
  @(0.3673045729500757,0.07342274756952349,0.6126623398149604,0.2355572969155837,0,0.00250353495305797,0,0.9532090292353743,0.1118950956651048,0.4240039180121471,0,0.3318957728824614,2.160070240037701,4.871888618469143),
  @(0.3483554653338956,0.07003953440113264,0.6096426966952606,0.2348225487910653,0,0.002505267865407249,0,0.9569659062567162,0.1118308127153433,0.4026866512202218,0,0.3252329329025585,2.172928595430559,4.884314632267916),
  @(0.3406654421252426,0.06866778465962398,0.6084838905982508,0.2345499226642396,0,0.002505996449479116,0,0.9585986958208501,0.1118166014571322,0.3940372884579801,0,0.3225581996136881,2.178326242399988,4.889928380854514),
  @(0.3393904584656582,0.06844042702576303,0.6082958077330858,0.2345062726014646,0,0.002506118785841227,0,0.9588759709598911,0.1118149659193612,0.3926033521278214,0,0.3221165087277882,2.179232046656832,4.89089376953163),
  @(0.3482516252923347,0.07002100635290276,0.6096267780410187,0.2348187635319761,0,0.002505277600542533,0,0.956987514258401,0.1118305725145028,0.4025698500761052,0,0.3251966966067528,2.173000751185796,4.88438811375002),
  @(0.3876328465865697,0.07705653292683223,0.6161293527958804,0.2364327938593753,0,0.002501762182593725,0,0.9495480501507707,0.1120046199729927,0.446878701223909,0,0.3391433806793458,2.146888385156171,4.860508271501743),
  @(0.4658373104905138,0.09106929496840621,0.6311348456373764,0.2404404121888746,0,0.002495570396622599,0,0.9381695457801769,0.1127232255134558,0.5349257820606965,0,0.3677579475282471,2.100657944102593,4.831017988131919),
  @(0.5238764720306506,0.101495578640808,0.6435289844991416,0.2438970845843471,0,0.002491444954061106,0,0.9317651975503551,0.11348080372629,0.6003065162577172,0,0.3895476517071401,2.069722474541059,4.819913815439861),
  @(0.5504038551531494,0.1062671661986769,0.6494639740931234,0.2455806071043014,0,0.002489659272779149,0,0.9292760082179896,0.1138752365736693,0.6301979069277195,0,0.3996259693558386,2.056306359458986,4.817154611844614),
  @(0.5604667160975794,0.1080781263959807,0.6517539745894965,0.2462340512604158,0,0.002488996099205543,0,0.9283943867158158,0.1140317505960908,0.641538121509484,0,0.4034661192085807,2.051320430715837,4.816439214513906),
  @(0.558298727361489,0.1076879238128186,0.651258891802911,0.2460926124886385,0,0.002489138347095851,0,0.9285815478102535,0.1139977246434967,0.6390948774246965,0,0.4026380224629236,2.05239003764239,4.816578638034997),
  @(0.5512313844081973,0.1064160738413307,0.6496515218486252,0.2456340473038239,0,0.002489604452436291,0,0.9292022546380494,0.1138879698267417,0.6311304552419017,0,0.3999414265995398,2.055894270645606,4.817089153816852),
  @(0.5469047000403577,0.1056375566662382,0.648672498417028,0.2453552363411404,0,0.002489891649633811,0,0.9295903964776429,0.1138216727519392,0.6262547362792361,0,0.3982927663859073,2.058053017834116,4.817444759475336),
  @(0.5221453219212719,0.1011843156275063,0.6431470804361368,0.2437892938947286,0,0.002491563478680398,0,0.9319364061327917,0.1134560278854693,0.5983560077596906,0,0.3888923371036981,2.070612465311196,4.820140246410261),
  @(0.5069879468772456,0.09845968959612605,0.6398333452606835,0.2428570566499744,0,0.002492612356194821,0,0.9334842382641924,0.1132444652609053,0.5812789527986126,0,0.3831678914335015,2.078485479051349,4.822380802350153),
  @(0.4982816265954852,0.09689525360340667,0.637955319935827,0.2423313148078066,0,0.002493224212244304,0,0.9344144372047438,0.1131274687517916,0.5714707705773208,0,0.379890983416054,2.08307562623007,4.823885266984689),
  @(0.4953358572646493,0.09636602721930387,0.6373242569684407,0.2421551047714843,0,0.002493432849785654,0,0.9347362445597156,0.1130886613797912,0.5681523230300058,0,0.3787841703628985,2.08464038367617,4.824431713914208),
  @(0.5086002573719952,0.09874945195279849,0.6401832065485564,0.2429552129937278,0,0.00249249981514684,0,0.9333153367310345,0.113266501280826,0.583095379897685,0,0.383775650827836,2.077640987113899,4.822119962509646),
  @(0.5533067626928982,0.1067895371376437,0.6501224916880233,0.2457683069804659,0,0.002489467192794526,0,0.9290182831307661,0.1139200134943437,0.6334692332028169,0,0.4007328400627586,2.054862428386258,4.816930262906709),
  @(0.5826269309528413,0.1120678590064017,0.6568663327523154,0.2476996642235534,0,0.002487561093309773,0,0.9265653439000801,0.1143887969968205,0.6665136001596181,0,0.4119534620831047,2.040525949173972,4.815458729672343),
  @(0.5669690381992893,0.1092485713824942,0.6532443758153761,0.2466603811685104,0,0.002488571490178448,0,0.9278420067559452,0.1141347888278332,0.6488661839923964,0,0.405952223711104,2.048127199448317,4.816068463249479),
  @(0.5078713072900598,0.09861844419566523,0.6400249498840651,0.2429108046933592,0,0.002492550667250816,0,0.9333915714866805,0.1132565243620505,0.5822741431297231,0,0.3835008387111998,2.078022583038248,4.822237214296962),
  @(0.4445776628689089,0.08725540140075339,0.6268346223550907,0.2392661875532873,0,0.002497170733524218,0,0.9409042220791406,0.1124884719984180,0.5109841466398564,0,0.3598819896724024,2.112632166146165,4.837140510911041),
)

$startRow = 2
$startCol = 2  # column B

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $startCol + $c).Value = $rowVals[$c]
    }
}
